$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D,E,F,G (open/close/high/low), H (shares_outstanding), I (fixed_ticker)
# for each data row (rows 2-44, row 32 already correct/unchanged).

$ws.Range("D2").Value = 6.920000076293945
$ws.Range("E2").Value = 7.130000114440918
$ws.Range("F2").Value = 7.659999847412109
$ws.Range("G2").Value = 6.449999809265137
$ws.Range("H2").Value = 136895573
$ws.Range("I2").Value = "LSCC"

$ws.Range("D3").Value = 6.340000152587891
$ws.Range("E3").Value = 5.929999828338623
$ws.Range("F3").Value = 6.579999923706055
$ws.Range("G3").Value = 5.900000095367432
$ws.Range("H3").Value = 136895573
$ws.Range("I3").Value = "LSCC"

$ws.Range("D4").Value = 5.929999828338623
$ws.Range("E4").Value = 4.920000076293945
$ws.Range("F4").Value = 6.139999866485596
$ws.Range("G4").Value = 4.510000228881836
$ws.Range("H4").Value = 136895573
$ws.Range("I4").Value = "LSCC"

$ws.Range("D5").Value = 3.829999923706055
$ws.Range("E5").Value = 4.579999923706055
$ws.Range("F5").Value = 4.880000114440918
$ws.Range("G5").Value = 3.569999933242798
$ws.Range("H5").Value = 136895573
$ws.Range("I5").Value = "LSCC"

$ws.Range("D6").Value = 6.400000095367432
$ws.Range("E6").Value = 4.860000133514404
$ws.Range("F6").Value = 6.400000095367432
$ws.Range("G6").Value = 4.019999980926514
$ws.Range("H6").Value = 136895573
$ws.Range("I6").Value = "LSCC"

$ws.Range("D7").Value = 5.619999885559082
$ws.Range("E7").Value = 5.570000171661377
$ws.Range("F7").Value = 6.46999979019165
$ws.Range("G7").Value = 5.139999866485596
$ws.Range("H7").Value = 136895573
$ws.Range("I7").Value = "LSCC"

$ws.Range("D8").Value = 5.360000133514404
$ws.Range("E8").Value = 6.010000228881836
$ws.Range("F8").Value = 6.300000190734863
$ws.Range("G8").Value = 5.210000038146973
$ws.Range("H8").Value = 136895573
$ws.Range("I8").Value = "LSCC"

$ws.Range("D9").Value = 6.449999809265137
$ws.Range("E9").Value = 6.070000171661377
$ws.Range("F9").Value = 6.650000095367432
$ws.Range("G9").Value = 5.909999847412109
$ws.Range("H9").Value = 136895573
$ws.Range("I9").Value = "LSCC"

$ws.Range("D10").Value = 7.389999866485596
$ws.Range("E10").Value = 7.190000057220459
$ws.Range("F10").Value = 7.420000076293945
$ws.Range("G10").Value = 6.929999828338623
$ws.Range("H10").Value = 136895573
$ws.Range("I10").Value = "LSCC"

$ws.Range("D11").Value = 6.940000057220459
$ws.Range("E11").Value = 6.860000133514404
$ws.Range("F11").Value = 7.010000228881836
$ws.Range("G11").Value = 6.659999847412109
$ws.Range("H11").Value = 136895573
$ws.Range("I11").Value = "LSCC"

$ws.Range("D12").Value = 6.699999809265137
$ws.Range("E12").Value = 6.960000038146973
$ws.Range("F12").Value = 6.989999771118164
$ws.Range("G12").Value = 6.510000228881836
$ws.Range("H12").Value = 136895573
$ws.Range("I12").Value = "LSCC"

$ws.Range("D13").Value = 5.25
$ws.Range("E13").Value = 5.849999904632568
$ws.Range("F13").Value = 6.010000228881836
$ws.Range("G13").Value = 5.139999866485596
$ws.Range("H13").Value = 136895573
$ws.Range("I13").Value = "LSCC"

$ws.Range("D14").Value = 5.78000020980835
$ws.Range("E14").Value = 6.510000228881836
$ws.Range("F14").Value = 6.71999979019165
$ws.Range("G14").Value = 5.78000020980835
$ws.Range("H14").Value = 136895573
$ws.Range("I14").Value = "LSCC"

$ws.Range("D15").Value = 5.570000171661377
$ws.Range("E15").Value = 5.420000076293945
$ws.Range("F15").Value = 5.679999828338623
$ws.Range("G15").Value = 5.050000190734863
$ws.Range("H15").Value = 136895573
$ws.Range("I15").Value = "LSCC"

$ws.Range("D16").Value = 6.510000228881836
$ws.Range("E16").Value = 7.690000057220459
$ws.Range("F16").Value = 7.980000019073486
$ws.Range("G16").Value = 6.28000020980835
$ws.Range("H16").Value = 136895573
$ws.Range("I16").Value = "LSCC"

$ws.Range("D17").Value = 8.0600004196167
$ws.Range("E17").Value = 6.010000228881836
$ws.Range("F17").Value = 8.899999618530273
$ws.Range("G17").Value = 5.579999923706055
$ws.Range("H17").Value = 136895573
$ws.Range("I17").Value = "LSCC"

$ws.Range("D18").Value = 6.860000133514404
$ws.Range("E18").Value = 7.800000190734863
$ws.Range("F18").Value = 7.929999828338623
$ws.Range("G18").Value = 6.710000038146973
$ws.Range("H18").Value = 136895573
$ws.Range("I18").Value = "LSCC"

$ws.Range("D19").Value = 12.06999969482422
$ws.Range("E19").Value = 12.94999980926514
$ws.Range("F19").Value = 13.57999992370606
$ws.Range("G19").Value = 11.1899995803833
$ws.Range("H19").Value = 136895573
$ws.Range("I19").Value = "LSCC"

$ws.Range("D20").Value = 15.13000011444092
$ws.Range("E20").Value = 19.34000015258789
$ws.Range("F20").Value = 20.25
$ws.Range("G20").Value = 14.47000026702881
$ws.Range("H20").Value = 136895573
$ws.Range("I20").Value = "LSCC"

$ws.Range("D21").Value = 18.42000007629395
$ws.Range("E21").Value = 19.59000015258789
$ws.Range("F21").Value = 20.09000015258789
$ws.Range("G21").Value = 17.05999946594238
$ws.Range("H21").Value = 136895573
$ws.Range("I21").Value = "LSCC"

$ws.Range("D22").Value = 19.45999908447266
$ws.Range("E22").Value = 18.60000038146973
$ws.Range("F22").Value = 24.20000076293945
$ws.Range("G22").Value = 18.45000076293945
$ws.Range("H22").Value = 136895573
$ws.Range("I22").Value = "LSCC"

$ws.Range("D23").Value = 17.03000068664551
$ws.Range("E23").Value = 22.51000022888184
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 16.42000007629395
$ws.Range("H23").Value = 136895573
$ws.Range("I23").Value = "LSCC"

$ws.Range("D24").Value = 28.21999931335449
$ws.Range("E24").Value = 31.09000015258789
$ws.Range("F24").Value = 31.79000091552734
$ws.Range("G24").Value = 27.81999969482422
$ws.Range("H24").Value = 136895573
$ws.Range("I24").Value = "LSCC"

$ws.Range("D25").Value = 29.46999931335449
$ws.Range("E25").Value = 34.90000152587891
$ws.Range("F25").Value = 35.90000152587891
$ws.Range("G25").Value = 28.05999946594238
$ws.Range("H25").Value = 136895573
$ws.Range("I25").Value = "LSCC"

$ws.Range("D26").Value = 46.04999923706055
$ws.Range("E26").Value = 40.11000061035156
$ws.Range("F26").Value = 46.88999938964844
$ws.Range("G26").Value = 37.38000106811523
$ws.Range("H26").Value = 136895573
$ws.Range("I26").Value = "LSCC"

$ws.Range("D27").Value = 45.84000015258789
$ws.Range("E27").Value = 50.31000137329102
$ws.Range("F27").Value = 58.38000106811523
$ws.Range("G27").Value = 45.84000015258789
$ws.Range("H27").Value = 136895573
$ws.Range("I27").Value = "LSCC"

$ws.Range("D28").Value = 56.04000091552734
$ws.Range("E28").Value = 56.75
$ws.Range("F28").Value = 56.84000015258789
$ws.Range("G28").Value = 49.04999923706055
$ws.Range("H28").Value = 136895573
$ws.Range("I28").Value = "LSCC"

$ws.Range("D29").Value = 64.94000244140625
$ws.Range("E29").Value = 69.44000244140625
$ws.Range("F29").Value = 70.94000244140625
$ws.Range("G29").Value = 62.56999969482422
$ws.Range("H29").Value = 136895573
$ws.Range("I29").Value = "LSCC"

$ws.Range("D30").Value = 77.54000091552734
$ws.Range("E30").Value = 55.22000122070312
$ws.Range("F30").Value = 79.22000122070312
$ws.Range("G30").Value = 47.04000091552734
$ws.Range("H30").Value = 136895573
$ws.Range("I30").Value = "LSCC"

$ws.Range("D31").Value = 61.29000091552734
$ws.Range("E31").Value = 48.04000091552734
$ws.Range("F31").Value = 62.06999969482422
$ws.Range("G31").Value = 46.47000122070312
$ws.Range("H31").Value = 136895573
$ws.Range("I31").Value = "LSCC"

$ws.Range("D33").Value = 49.68000030517578
$ws.Range("E33").Value = 48.5099983215332
$ws.Range("F33").Value = 56.66999816894531
$ws.Range("G33").Value = 44.52999877929688
$ws.Range("H33").Value = 136895573
$ws.Range("I33").Value = "LSCC"

$ws.Range("D34").Value = 66.87000274658203
$ws.Range("E34").Value = 75.79000091552734
$ws.Range("F34").Value = 77.66000366210938
$ws.Range("G34").Value = 63.4900016784668
$ws.Range("H34").Value = 136895573
$ws.Range("I34").Value = "LSCC"

$ws.Range("D35").Value = 94.80999755859376
$ws.Range("E35").Value = 79.69999694824219
$ws.Range("F35").Value = 96.16999816894533
$ws.Range("G35").Value = 75.86000061035156
$ws.Range("H35").Value = 136895573
$ws.Range("I35").Value = "LSCC"

$ws.Range("D36").Value = 96.31999969482422
$ws.Range("E36").Value = 90.94000244140624
$ws.Range("F36").Value = 98.18000030517578
$ws.Range("G36").Value = 84.69999694824219
$ws.Range("H36").Value = 136895573
$ws.Range("I36").Value = "LSCC"

$ws.Range("D37").Value = 86.33999633789062
$ws.Range("E37").Value = 55.61000061035156
$ws.Range("F37").Value = 87.63999938964844
$ws.Range("G37").Value = 53.7400016784668
$ws.Range("H37").Value = 136895573
$ws.Range("I37").Value = "LSCC"

$ws.Range("D38").Value = 68.16000366210938
$ws.Range("E38").Value = 60.86000061035156
$ws.Range("F38").Value = 71.05000305175781
$ws.Range("G38").Value = 59.34999847412109
$ws.Range("H38").Value = 136895573
$ws.Range("I38").Value = "LSCC"

$ws.Range("D39").Value = 78.25
$ws.Range("E39").Value = 68.5999984741211
$ws.Range("F39").Value = 80.12999725341797
$ws.Range("G39").Value = 65.38999938964844
$ws.Range("H39").Value = 136895573
$ws.Range("I39").Value = "LSCC"

$ws.Range("D40").Value = 57.88999938964844
$ws.Range("E40").Value = 53
$ws.Range("F40").Value = 63.20000076293945
$ws.Range("G40").Value = 47.58000183105469
$ws.Range("H40").Value = 136895573
$ws.Range("I40").Value = "LSCC"

$ws.Range("D41").Value = 53.15999984741211
$ws.Range("E41").Value = 50.65999984741211
$ws.Range("F41").Value = 55.83000183105469
$ws.Range("G41").Value = 48.22999954223633
$ws.Range("H41").Value = 136895573
$ws.Range("I41").Value = "LSCC"

$ws.Range("D42").Value = 57.2599983215332
$ws.Range("E42").Value = 57.02000045776367
$ws.Range("F42").Value = 61.88000106811523
$ws.Range("G42").Value = 52.61999893188477
$ws.Range("H42").Value = 136895573
$ws.Range("I42").Value = "LSCC"

$ws.Range("D43").Value = 52.22999954223633
$ws.Range("E43").Value = 48.93000030517578
$ws.Range("F43").Value = 53.18000030517578
$ws.Range("G43").Value = 34.68999862670898
$ws.Range("H43").Value = 136895573
$ws.Range("I43").Value = "LSCC"

$ws.Range("D44").Value = 48.22999954223633
$ws.Range("E44").Value = 49.83000183105469
$ws.Range("F44").Value = 56.75
$ws.Range("G44").Value = 48.09999847412109
$ws.Range("H44").Value = 136895573
$ws.Range("I44").Value = "LSCC"
